# Slide 3, "TextBox 4" shape: clarify that Thunder_data_per_station.ipynb is
# the (in-progress) notebook being used to produce thunder_counts.py /
# thunder_averages.py.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item("TextBox 4")

# Prepend two new paragraphs ahead of the existing "thunder_counts.py" /
# "thunder_averages.py" lines, matching how the author typed them in (a
# line naming the notebook, then a tab-indented status line).
$sh.TextFrame.TextRange.InsertBefore("Thunder_data_per_station.ipynb`r`tcurrently being converted:`r")

# The textbox grew taller (spAutoFit) to fit the extra lines and was nudged
# left so it stays centered over the same content.
$emuPerPt = 914400 / 72
$sh.Left = 2957869 / $emuPerPt
